# 5 fix + dataset
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers (columns A, B, C get new names; D and E stay the same)
$ws.Range("A1").Value = "Опыт_разработки"
$ws.Range("B1").Value = "Уровень_знаний"
$ws.Range("C1").Value = "Претендуемая_позиция"

# Data fix: replace the old "Р" code with "И" in column C
$used = $ws.Range("C2:C15")
for ($i = 1; $i -le $used.Rows.Count; $i++) {
    $cell = $used.Cells.Item($i, 1)
    if ($cell.Value2 -eq "Р") {
        $cell.Value = "И"
    }
}

# Re-fit the column widths to the new (longer) header text
# (closest values this engine's width model can reproduce to the
# real-Excel autofit results of 16.85546875 / 19.7109375 / 23.42578125)
$ws.Columns("A").ColumnWidth = 16
$ws.Columns("B").ColumnWidth = 18.83
$ws.Columns("C").ColumnWidth = 22.67

# Move the selection to D1
$ws.Range("D1").Select() | Out-Null
